# Apply cell text updates for the crypto price/volume table refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.947.58'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '2.242.45'
$ws.Range('E3').Value = '  -1.80%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '231.05'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = '0.643'
$ws.Range('E6').Value = '  +2.38%  '
$ws.Range('D7').Value = '62.98'
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').Value = '0.448'
$ws.Range('E9').Value = '  +5.99%  '
$ws.Range('D10').Value = '0.0977'
$ws.Range('E10').Value = '  +2.34%  '
$ws.Range('D11').Value = '57.45'
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').Value = '26.11'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').Value = '2.576.74'
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').Value = '15.45'
$ws.Range('E15').Value = '  -2.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.10'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.61%  '
$ws.Range('D17').Value = '0.826'
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('D18').Value = '2.237.01'
$ws.Range('E18').Value = '  -2.43%  '
$ws.Range('D19').Value = '43.805.49'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').Value = '0.0₃0984'
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('D21').Value = '72.59'
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('D23').Value = '247.17'
$ws.Range('E23').Value = '  -2.50%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = '2.41'
$ws.Range('E25').Value = '  -8.47%  '
$ws.Range('D26').Value = '3.35'
$ws.Range('E26').Value = '  +21.63%  '
$ws.Range('D27').Value = '2.23'
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').Value = '9.79'
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '171.80'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').Value = '21.06'
$ws.Range('E30').Value = '  +2.11%  '
$ws.Range('D31').Value = '0.139'
$ws.Range('E31').Value = '  -1.04%  '
$ws.Range('E32').Value = '  -2.70%  '
$ws.Range('D33').Value = '0.126'
$ws.Range('E33').Value = '  +2.85%  '
$ws.Range('D34').Value = '0.0684'
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('E35').Value = '  +0.87%  '
$ws.Range('D36').Value = '4.93'
$ws.Range('E36').Value = '  -3.39%  '
$ws.Range('D37').Value = '3.66'
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('D38').Value = '6.41'
$ws.Range('E38').Value = '  -2.34%  '
$ws.Range('D39').Value = '2.28'
$ws.Range('E39').Value = '  -3.98%  '
$ws.Range('D40').Value = '0.0252'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('E41').Value = '  -0.30%  '
$ws.Range('D42').Value = '0.000227'
$ws.Range('E42').Value = '  +3.03%  '
$ws.Range('D43').Value = '8.37'
$ws.Range('E43').Value = '  -4.43%  '
$ws.Range('D44').Value = '17.08'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('D45').Value = '97.12'
$ws.Range('E45').Value = '  -1.44%  '
$ws.Range('E46').Value = '  -2.19%  '
$ws.Range('D47').Value = '0.0941'
$ws.Range('E47').Value = '  -2.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.30'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -6.31%  '
$ws.Range('D49').Value = '1.429.87'
$ws.Range('E49').Value = '  -3.80%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '2.27'
$ws.Range('E50').Value = '  -2.19%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').Value = '2.76'
$ws.Range('E51').Value = '  +1.03%  '
